$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E10: "the number of PubMed articles about this Species", with "PubMed"
# as a hyperlink to PubMed (rich text: only "PubMed" is blue/linked) ---
$cell = $ws.Range("E10")

# Add the hyperlink first so the later Value/Characters calls are the ones
# that stick (Hyperlinks.Add overwrites the cell text with its display text).
$ws.Hyperlinks.Add($cell, "https://pubmed.ncbi.nlm.nih.gov/", "", "PubMed", "PubMed")

$cell.Value = "the number of PubMed articles about this Species"
$chars = $cell.Characters(15, 6)
$chars.Font.Name = "Arial"
$chars.Font.Size = 10
$chars.Font.Color = 16711680
$chars.Font.Underline = 0

# --- Selection moved from C21 to C15 ---
$ws.Range("C15").Select()

# --- Column D width: 16.15 -> 16.14 (closest value reachable through the
# character-width rounding used by this host) ---
$ws.Columns.Item(4).ColumnWidth = 15.25
